# Auto-generated update of cryptos list (mirrors the upstream GitHub Actions scrape commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.371.88"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.842.60"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D5").Value = "'238.78"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'0.6306"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "'0.2925"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "'24.36"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "1.878.63"
$ws.Range("E12").Value = "  -5.36%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "'0.00001028"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "'82.87"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "2.119.78"
$ws.Range("E17").Value = "  -6.38%  "
$ws.Range("D18").Value = "'6.148"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "29.412.48"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'227.58"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'7.440"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("D25").Value = "'156.92"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'8.355"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "'17.60"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "'1.455"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").Value = "'1.276"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'1.829"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").Value = "'1.155"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'0.7061"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "1.240.68"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'2.764"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").Value = "'6.313"
$ws.Range("E41").Value = "  +3.62%  "
$ws.Range("D42").Value = "'0.8999"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'0.9994"
$ws.Range("D44").Value = "'101.89"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").Value = "'65.58"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.059"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").Value = "'0.3999"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.665"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.871"
$ws.Range("E49").Value = "  -3.42%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1118"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05708"
$ws.Range("E51").Value = "  -0.64%  "
